# Scheduled market-data refresh: updates currentAveragePrice* / Leve price & profit
# columns (H:N) for the affected leve rows on each job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Cells.Item(8, 8).Value = 37.1
$ws.Cells.Item(8, 9).Value = 23.714285
$ws.Cells.Item(8, 10).Value = 68.333336
$ws.Cells.Item(8, 11).Value = 71.142855
$ws.Cells.Item(8, 12).Value = 205.000008
$ws.Cells.Item(8, 13).Value = 67.857145
$ws.Cells.Item(8, 14).Value = -483.000008
# Row 20
$ws.Cells.Item(20, 8).Value = 721.55554
$ws.Cells.Item(20, 9).Value = 761.625
$ws.Cells.Item(20, 10).Value = 401
$ws.Cells.Item(20, 11).Value = 761.625
$ws.Cells.Item(20, 12).Value = 401
$ws.Cells.Item(20, 13).Value = -531.625
$ws.Cells.Item(20, 14).Value = -861
# Row 35
$ws.Cells.Item(35, 8).Value = 721.55554
$ws.Cells.Item(35, 9).Value = 761.625
$ws.Cells.Item(35, 10).Value = 401
$ws.Cells.Item(35, 11).Value = 761.625
$ws.Cells.Item(35, 12).Value = 401
$ws.Cells.Item(35, 13).Value = -382.625
$ws.Cells.Item(35, 14).Value = -1159
# Row 42
$ws.Cells.Item(42, 8).Value = 122.125
$ws.Cells.Item(42, 9).Value = 129.5
$ws.Cells.Item(42, 11).Value = 388.5
$ws.Cells.Item(42, 13).Value = -158.5
# Row 70
$ws.Cells.Item(70, 8).Value = 4213.4546
$ws.Cells.Item(70, 9).Value = 1000
$ws.Cells.Item(70, 11).Value = 3000
$ws.Cells.Item(70, 13).Value = -2730
# Row 73
$ws.Cells.Item(73, 8).Value = 4213.4546
$ws.Cells.Item(73, 9).Value = 1000
$ws.Cells.Item(73, 11).Value = 3000
$ws.Cells.Item(73, 13).Value = -2064
# Row 74
$ws.Cells.Item(74, 8).Value = 86800
$ws.Cells.Item(74, 9).Value = 4160
$ws.Cells.Item(74, 11).Value = 4160
$ws.Cells.Item(74, 13).Value = -3224
# Row 76
$ws.Cells.Item(76, 8).Value = 3333.3333
$ws.Cells.Item(76, 9).Value = 4000
$ws.Cells.Item(76, 10).Value = 3000
$ws.Cells.Item(76, 11).Value = 4000
$ws.Cells.Item(76, 12).Value = 3000
$ws.Cells.Item(76, 13).Value = -3685
$ws.Cells.Item(76, 14).Value = -3630
# Row 77
$ws.Cells.Item(77, 8).Value = 86800
$ws.Cells.Item(77, 9).Value = 4160
$ws.Cells.Item(77, 11).Value = 20800
$ws.Cells.Item(77, 13).Value = -16120
# Row 79
$ws.Cells.Item(79, 8).Value = 3333.3333
$ws.Cells.Item(79, 9).Value = 4000
$ws.Cells.Item(79, 10).Value = 3000
$ws.Cells.Item(79, 11).Value = 4000
$ws.Cells.Item(79, 12).Value = 3000
$ws.Cells.Item(79, 13).Value = -2908
$ws.Cells.Item(79, 14).Value = -5184
# Row 120
$ws.Cells.Item(120, 8).Value = 27500
$ws.Cells.Item(120, 10).Value = 27500
$ws.Cells.Item(120, 12).Value = 27500
$ws.Cells.Item(120, 14).Value = -37176
# Row 125
$ws.Cells.Item(125, 9).Value = 2390
$ws.Cells.Item(125, 10).Value = 8383.5
$ws.Cells.Item(125, 11).Value = 21510
$ws.Cells.Item(125, 12).Value = 75451.5
$ws.Cells.Item(125, 13).Value = -19050
$ws.Cells.Item(125, 14).Value = -80371.5

$ws = $wb.Worksheets.Item("ARM")
# Row 8
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 9).Value = 0
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 0
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 13).Value = $null
$ws.Cells.Item(8, 14).Value = $null
# Row 25
$ws.Cells.Item(25, 8).Value = 1578
$ws.Cells.Item(25, 9).Value = 722.5
$ws.Cells.Item(25, 11).Value = 722.5
$ws.Cells.Item(25, 13).Value = -320.5
# Row 35
$ws.Cells.Item(35, 8).Value = 3406
$ws.Cells.Item(35, 9).Value = 2115.5
$ws.Cells.Item(35, 10).Value = 4266.3335
$ws.Cells.Item(35, 11).Value = 2115.5
$ws.Cells.Item(35, 12).Value = 4266.3335
$ws.Cells.Item(35, 13).Value = -1709.5
$ws.Cells.Item(35, 14).Value = -5078.3335
# Row 88
$ws.Cells.Item(88, 8).Value = 2947.7856
$ws.Cells.Item(88, 9).Value = 1317
$ws.Cells.Item(88, 11).Value = 1317
$ws.Cells.Item(88, 13).Value = -911
# Row 91
$ws.Cells.Item(91, 8).Value = 2947.7856
$ws.Cells.Item(91, 9).Value = 1317
$ws.Cells.Item(91, 11).Value = 1317
$ws.Cells.Item(91, 13).Value = 87

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 14).Value = $null
# Row 37
$ws.Cells.Item(37, 8).Value = 955.5
$ws.Cells.Item(37, 9).Value = 892
$ws.Cells.Item(37, 10).Value = 1400
$ws.Cells.Item(37, 11).Value = 892
$ws.Cells.Item(37, 12).Value = 1400
$ws.Cells.Item(37, 13).Value = -755
$ws.Cells.Item(37, 14).Value = -1674
# Row 46
$ws.Cells.Item(46, 8).Value = 10000
$ws.Cells.Item(46, 10).Value = 10000
$ws.Cells.Item(46, 12).Value = 10000
$ws.Cells.Item(46, 14).Value = -10596
# Row 105
$ws.Cells.Item(105, 8).Value = 746.5
$ws.Cells.Item(105, 9).Value = 785
$ws.Cells.Item(105, 10).Value = 708
$ws.Cells.Item(105, 11).Value = 785
$ws.Cells.Item(105, 12).Value = 708
$ws.Cells.Item(105, 13).Value = 962
$ws.Cells.Item(105, 14).Value = -4202

$ws = $wb.Worksheets.Item("CRP")
# Row 42
$ws.Cells.Item(42, 8).Value = 0
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(42, 13).Value = $null
# Row 50
$ws.Cells.Item(50, 8).Value = 29998.75
$ws.Cells.Item(50, 10).Value = 29998.75
$ws.Cells.Item(50, 12).Value = 29998.75
$ws.Cells.Item(50, 14).Value = -31248.75
# Row 62
$ws.Cells.Item(62, 8).Value = 1900
$ws.Cells.Item(62, 9).Value = 1900
$ws.Cells.Item(62, 11).Value = 1900
$ws.Cells.Item(62, 13).Value = -1276
# Row 65
$ws.Cells.Item(65, 8).Value = 1900
$ws.Cells.Item(65, 9).Value = 1900
$ws.Cells.Item(65, 11).Value = 9500
$ws.Cells.Item(65, 13).Value = -6380
# Row 107
$ws.Cells.Item(107, 8).Value = 1225.2858
$ws.Cells.Item(107, 9).Value = 1305.8334
$ws.Cells.Item(107, 11).Value = 1305.8334
$ws.Cells.Item(107, 13).Value = 614.1666

$ws = $wb.Worksheets.Item("CUL")
# Row 6
$ws.Cells.Item(6, 8).Value = 15251.75
$ws.Cells.Item(6, 9).Value = 504
$ws.Cells.Item(6, 11).Value = 1512
$ws.Cells.Item(6, 13).Value = -1399
# Row 10
$ws.Cells.Item(10, 8).Value = 350.91666
$ws.Cells.Item(10, 9).Value = 19.363636
$ws.Cells.Item(10, 11).Value = 58.090908
$ws.Cells.Item(10, 13).Value = 80.909092
# Row 11
$ws.Cells.Item(11, 8).Value = 1037.2222
$ws.Cells.Item(11, 9).Value = 1065.4
$ws.Cells.Item(11, 11).Value = 3196.2
$ws.Cells.Item(11, 13).Value = -3056.2
# Row 36
$ws.Cells.Item(36, 8).Value = 540.6667
$ws.Cells.Item(36, 9).Value = 361
$ws.Cells.Item(36, 11).Value = 1083
$ws.Cells.Item(36, 13).Value = -914
# Row 55
$ws.Cells.Item(55, 8).Value = 1004
$ws.Cells.Item(55, 9).Value = 1004
$ws.Cells.Item(55, 11).Value = 3012
$ws.Cells.Item(55, 13).Value = -2835
# Row 95
$ws.Cells.Item(95, 8).Value = 3027
$ws.Cells.Item(95, 10).Value = 3027
$ws.Cells.Item(95, 12).Value = 9081
$ws.Cells.Item(95, 14).Value = -13199
# Row 113
$ws.Cells.Item(113, 8).Value = 950.1667
$ws.Cells.Item(113, 10).Value = 924.5
$ws.Cells.Item(113, 12).Value = 2773.5
$ws.Cells.Item(113, 14).Value = -7113.5
# Row 138
$ws.Cells.Item(138, 8).Value = 6705.294
$ws.Cells.Item(138, 9).Value = 3996.3333
$ws.Cells.Item(138, 10).Value = 7285.7856
$ws.Cells.Item(138, 11).Value = 11988.9999
$ws.Cells.Item(138, 12).Value = 21857.3568
$ws.Cells.Item(138, 13).Value = -6848.999899999999
$ws.Cells.Item(138, 14).Value = -32137.3568

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 14).Value = $null
# Row 22
$ws.Cells.Item(22, 8).Value = 141339.33
$ws.Cells.Item(22, 9).Value = 4000
$ws.Cells.Item(22, 11).Value = 4000
$ws.Cells.Item(22, 13).Value = -3471
# Row 43
$ws.Cells.Item(43, 8).Value = 32154.4
$ws.Cells.Item(43, 10).Value = 32154.4
$ws.Cells.Item(43, 12).Value = 32154.4
$ws.Cells.Item(43, 14).Value = -32456.4
# Row 52
$ws.Cells.Item(52, 8).Value = 0
$ws.Cells.Item(52, 10).Value = 0
$ws.Cells.Item(52, 12).Value = 0
$ws.Cells.Item(52, 14).Value = $null
# Row 122
$ws.Cells.Item(122, 8).Value = 9647236
$ws.Cells.Item(122, 9).Value = 12540408
$ws.Cells.Item(122, 11).Value = 37621224
$ws.Cells.Item(122, 13).Value = -37618774
# Row 126
$ws.Cells.Item(126, 8).Value = 2359.3333
$ws.Cells.Item(126, 9).Value = 2039
$ws.Cells.Item(126, 11).Value = 6117
$ws.Cells.Item(126, 13).Value = -3647

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Cells.Item(16, 8).Value = 2359.9285
$ws.Cells.Item(16, 9).Value = 1105.5555
$ws.Cells.Item(16, 10).Value = 4617.8
$ws.Cells.Item(16, 11).Value = 1105.5555
$ws.Cells.Item(16, 12).Value = 4617.8
$ws.Cells.Item(16, 13).Value = -935.5554999999999
$ws.Cells.Item(16, 14).Value = -4957.8
# Row 22
$ws.Cells.Item(22, 8).Value = 1122.7273
$ws.Cells.Item(22, 9).Value = 907.2857
$ws.Cells.Item(22, 10).Value = 1499.75
$ws.Cells.Item(22, 11).Value = 907.2857
$ws.Cells.Item(22, 12).Value = 1499.75
$ws.Cells.Item(22, 13).Value = -612.2857
$ws.Cells.Item(22, 14).Value = -2089.75
# Row 27
$ws.Cells.Item(27, 8).Value = 1122.7273
$ws.Cells.Item(27, 9).Value = 907.2857
$ws.Cells.Item(27, 10).Value = 1499.75
$ws.Cells.Item(27, 11).Value = 907.2857
$ws.Cells.Item(27, 12).Value = 1499.75
$ws.Cells.Item(27, 13).Value = -800.2857
$ws.Cells.Item(27, 14).Value = -1713.75
# Row 30
$ws.Cells.Item(30, 8).Value = 460.7
$ws.Cells.Item(30, 9).Value = 323.625
$ws.Cells.Item(30, 11).Value = 323.625
$ws.Cells.Item(30, 13).Value = -215.625
# Row 35
$ws.Cells.Item(35, 8).Value = 1692.3334
$ws.Cells.Item(35, 9).Value = 1692.3334
$ws.Cells.Item(35, 10).Value = 0
$ws.Cells.Item(35, 11).Value = 1692.3334
$ws.Cells.Item(35, 12).Value = 0
$ws.Cells.Item(35, 13).Value = -1356.3334
$ws.Cells.Item(35, 14).Value = $null
# Row 122
$ws.Cells.Item(122, 8).Value = 2986
$ws.Cells.Item(122, 9).Value = 2300.6
$ws.Cells.Item(122, 10).Value = 4699.5
$ws.Cells.Item(122, 11).Value = 6901.799999999999
$ws.Cells.Item(122, 12).Value = 14098.5
$ws.Cells.Item(122, 13).Value = -4451.799999999999
$ws.Cells.Item(122, 14).Value = -18998.5

$ws = $wb.Worksheets.Item("WVR")
# Row 21
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 13).Value = $null
# Row 35
$ws.Cells.Item(35, 8).Value = 0
$ws.Cells.Item(35, 9).Value = 0
$ws.Cells.Item(35, 11).Value = 0
$ws.Cells.Item(35, 13).Value = $null
# Row 114
$ws.Cells.Item(114, 8).Value = 49999.5
$ws.Cells.Item(114, 10).Value = 49999.5
$ws.Cells.Item(114, 12).Value = 49999.5
$ws.Cells.Item(114, 14).Value = -58677.5
# Row 126
$ws.Cells.Item(126, 8).Value = 2405.6428
$ws.Cells.Item(126, 9).Value = 2118.7
$ws.Cells.Item(126, 10).Value = 3123
$ws.Cells.Item(126, 11).Value = 6356.099999999999
$ws.Cells.Item(126, 12).Value = 9369
$ws.Cells.Item(126, 13).Value = -3886.099999999999
$ws.Cells.Item(126, 14).Value = -14309
